$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the JSON-ish parameter values to proper JSON (quote keys/values)
$ws.Range("C2").Value = '{"username":"SwatiChetty","password":"123456"}'
$ws.Range("D2").Value = '{"textToValidate":"Welcome Vivek!!"}'
$ws.Range("C3").Value = '{"username":"Swati","password":"123"}'
$ws.Range("D3").Value = '{"textToValidate":"Order Create Successfully"}'

# Widen column C independently of column D (target stored width ~46.45; the
# engine quantizes ColumnWidth to 1/6-character steps, so 45.6667 is the
# input that lands on the nearest reachable grid point, 46.5)
$ws.Columns("C").ColumnWidth = 45.666666666666664

# Move the active selection to C9
$ws.Range("C9").Select() | Out-Null
